$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '65.250.88'
$ws.Range("E2").Value = '  +2.28%  '

# Row 3
$ws.Range("D3").Value = '3.172.30'
$ws.Range("E3").Value = '  +4.01%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.36%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.70'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.08%  '

# Row 7
$ws.Range("E7").Value = '  -0.01%  '

# Row 8
$ws.Range("D8").Value = '3.170.35'
$ws.Range("E8").Value = '  +4.08%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.532'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.27%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.162'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.58%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.20'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.14%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.504'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.72%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000272'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +18.84%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.62'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.05%  '

# Row 15
$ws.Range("D15").Value = '3.694.40'
$ws.Range("E15").Value = '  +4.07%  '

# Row 16
$ws.Range("D16").Value = '65.369.92'
$ws.Range("E16").Value = '  +2.36%  '

# Row 19
$ws.Range("E19").Value = '  +1.38%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '513.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.11%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.90'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.68%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.728'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.81%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.42'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.52%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.85'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.47%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.22'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.36%  '

# Row 26
$ws.Range("E26").Value = '  +0.03%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.12'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +12.87%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.94'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.36%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.20'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +8.81%  '

# Row 30
$ws.Range("E30").Value = '  +16.38%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '27.91'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.39%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.10%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.20'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.36%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.32'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +11.76%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.62'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.02%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.86'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.29%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0906'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +11.48%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '479.28'
$ws.Range("D38").Style = "Normal"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.11'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +12.93%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0422'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.57%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.69'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.79%  '

# Row 42
$ws.Range("D42").Value = '3.086.96'
$ws.Range("E42").Value = '  +2.29%  '

# Row 43
$ws.Range("E43").Value = '  +3.10%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.48'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +11.89%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.286'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.25%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '29.41'
$ws.Range("D46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = '0.0₃0608'
$ws.Range("E47").Value = '  +18.92%  '

# Row 48
$ws.Range("E48").Value = '  -0.03%  '

# Row 49
$ws.Range("E49").Value = '  +1.55%  '

# Row 50
$ws.Range("E50").Value = '  +8.74%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.68'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.10%  '

# Row 17/18 swap (Polkadot <-> WrappedEther)
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.189.28'
$ws.Range("E17").Value = '  +4.50%  '

$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.19'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +6.52%  '

